$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands in the
# same tab order as the target workbook (Sheet1, Sheet2, SheetX_OnlyInFile1).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SheetX_OnlyInFile1"

$ws.Range("A1").Value = 123
$ws.Range("A2").Formula = "=A1 * 2"
